# Collapse the "(Socio)Technological rule or the Takeaway in this form: <br><tab>"
# run sequence into a single "Takeaway: " run, keeping the rest of the
# paragraph ("To achieve an effect in a situation apply this intervention ")
# untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 146 (top-level shape 4) -> GroupItems(3) is "Shape 149", the
# textbox holding the takeaway sentence.
$grp = $s.Shapes.Item(4)
$shp = $grp.GroupItems.Item(3)
$tr = $shp.TextFrame.TextRange

# The first 58 characters of the paragraph are:
#   "(Socio)" + "Technological rule" + " or the Takeaway in this form"
#   + ": " + "" + <line break> + "\t"
# Replacing that whole span with "Takeaway: " (using the formatting of
# the very first run) collapses all of those runs/br into one run and
# leaves the remainder of the paragraph untouched.
$target = $tr.Characters(1, 58)
$target.Text = "Takeaway: "
